$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49, shifting existing rows 49-121 down to 50-122.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record's data.
$ws.Cells.Item(49, 1).Value = 5
$ws.Cells.Item(49, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(49, 3).Value = "Maule"
$ws.Cells.Item(49, 4).Value = 44546
$ws.Cells.Item(49, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(49, 5).Value = 7
$ws.Cells.Item(49, 6).Value = 100112031
$ws.Cells.Item(49, 7).Value = "Poroto verde"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 300
$ws.Cells.Item(49, 11).Value = 16000
$ws.Cells.Item(49, 12).Value = 16000
$ws.Cells.Item(49, 13).Value = 16000
$ws.Cells.Item(49, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(49, 15).Value = "Región del Maule"
$ws.Cells.Item(49, 16).Value = 640
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = "Hortaliza"
